$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching style of existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in "Save" values for rows 2-8
$saveValues = @(0, 1, 1, 1, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
